# Filter - Study - Test Suit
# Change the "CasesTab" row label (cell A2 on the "startup" sheet) to
# "ParticipantsTab", and move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the cell value (was "CasesTab")
$ws.Range("A2").Value = "ParticipantsTab"

# Update the active selection shown in the sheet view
$ws.Activate()
$ws.Range("A2").Select()
